$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '28.405.59'
$ws.Cells.Item(2, 5).Value = '  +0.02%  '
$ws.Cells.Item(3, 4).Value = '1.565.65'
$ws.Cells.Item(3, 5).Value = '  -0.20%  '
$ws.Cells.Item(4, 5).Value = '  -0.13%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '210.81'
$ws.Cells.Item(5, 5).Value = '  -0.37%  '
$ws.Cells.Item(6, 5).Value = '  -0.27%  '
$ws.Cells.Item(7, 5).Value = '  -0.10%  '
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '44.35'
$ws.Cells.Item(8, 5).Value = '  -2.28%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '23.52'
$ws.Cells.Item(9, 5).Value = '  -1.88%  '
$ws.Cells.Item(10, 5).Value = '  -0.82%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.0586'
$ws.Cells.Item(11, 5).Value = '  -0.54%  '
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.0889'
$ws.Cells.Item(12, 5).Value = '  +0.00%  '
$ws.Cells.Item(13, 4).Value = '1.789.54'
$ws.Cells.Item(13, 5).Value = '  -0.19%  '
$ws.Cells.Item(14, 4).Value = '1.570.20'
$ws.Cells.Item(14, 5).Value = '  +0.12%  '
$ws.Cells.Item(15, 4).Value = '28.364.49'
$ws.Cells.Item(15, 5).Value = '  -0.16%  '
$ws.Cells.Item(16, 5).Value = '  -0.22%  '
$ws.Cells.Item(17, 5).Value = '  -1.35%  '
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '60.55'
$ws.Cells.Item(18, 5).Value = '  -2.57%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '228.96'
$ws.Cells.Item(19, 5).Value = '  +0.52%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '7.35'
$ws.Cells.Item(20, 5).Value = '  +0.47%  '
$ws.Cells.Item(21, 4).Value = '0.0₃0678'
$ws.Cells.Item(21, 5).Value = '  -1.46%  '
$ws.Cells.Item(22, 5).Value = '  -0.07%  '
$ws.Cells.Item(23, 5).Value = '  +1.51%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '8.86'
$ws.Cells.Item(24, 5).Value = '  -2.32%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '2.02'
$ws.Cells.Item(25, 5).Value = '  -4.43%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '150.19'
$ws.Cells.Item(26, 5).Value = '  -0.20%  '
$ws.Cells.Item(27, 5).Value = '  -0.52%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '0.104'
$ws.Cells.Item(28, 5).Value = '  +0.16%  '
$ws.Cells.Item(29, 2).Value = 'Cosmos'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '6.28'
$ws.Cells.Item(29, 5).Value = '  -2.32%  '
$ws.Cells.Item(30, 2).Value = 'BinanceUSD'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '1.00'
$ws.Cells.Item(30, 5).Value = '  -0.12%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '0.0471'
$ws.Cells.Item(31, 5).Value = '  -2.64%  '
$ws.Cells.Item(32, 5).Value = '  -2.71%  '
$ws.Cells.Item(33, 5).Value = '  -0.26%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '3.05'
$ws.Cells.Item(34, 5).Value = '  -0.94%  '
$ws.Cells.Item(35, 4).Value = '1.386.93'
$ws.Cells.Item(35, 5).Value = '  -0.49%  '
$ws.Cells.Item(36, 5).Value = '  +1.59%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '1.50'
$ws.Cells.Item(37, 5).Value = '  -2.45%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '2.34'
$ws.Cells.Item(38, 5).Value = '  -0.57%  '
$ws.Cells.Item(39, 5).Value = '  +1.46%  '
$ws.Cells.Item(40, 5).Value = '  -1.97%  '
$ws.Cells.Item(41, 2).Value = 'RenderToken'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '1.94'
$ws.Cells.Item(41, 5).Value = '  +2.67%  '
$ws.Cells.Item(42, 2).Value = 'ImmutableX'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '0.514'
$ws.Cells.Item(42, 5).Value = '  -2.92%  '
$ws.Cells.Item(43, 5).Value = '  -0.08%  '
$ws.Cells.Item(44, 5).Value = '  -0.32%  '
$ws.Cells.Item(45, 5).Value = '  -1.64%  '
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '5.38'
$ws.Cells.Item(46, 5).Value = '  -1.34%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '62.15'
$ws.Cells.Item(47, 5).Value = '  -0.53%  '
$ws.Cells.Item(48, 5).Value = '  -6.01%  '
$ws.Cells.Item(49, 4).Value = '1.701.86'
$ws.Cells.Item(49, 5).Value = '  -0.15%  '
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '85.59'
$ws.Cells.Item(50, 5).Value = '  -0.25%  '
$ws.Cells.Item(51, 2).Value = 'BitcoinSV'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '41.93'
$ws.Cells.Item(51, 5).Value = '  +10.45%  '
